# Horarios Línea 141 - actualización de datos (scrape 05:49:40)
# Agrega las nuevas filas relevadas a las hojas "LP1912" y "6203-6173",
# luego reordena cada tabla por "Hora_Llegada" (columna B) ascendente,
# y actualiza los textos de cabecera (última actualización / total filas).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Hoja "LP1912"
# ---------------------------------------------------------------
$wsLP = $wb.Worksheets.Item("LP1912")

# Filas existentes de datos: 6..42 (37 filas). Las nuevas se agregan a
# continuación (43..53) y luego se reordena todo el bloque 6..53.
$newRowsLP = @(
  ,@("05:49:40","05:51","14_ABASTO",2,"LP1912")
  ,@("05:49:40","06:06","16_SANTA ANA",17,"LP1912")
  ,@("05:49:40","07:04","23_HERNANDEZ",75,"LP1912")
  ,@("05:49:40","07:05","15_ABASTO",76,"LP1912")
  ,@("05:49:40","07:07","225_GOMEZ",78,"LP1912")
  ,@("05:49:40","07:21","26_HERNANDEZ",92,"LP1912")
  ,@("05:49:40","07:29","10_OLMOS",100,"LP1912")
  ,@("05:49:40","07:31","11_ETCHEVERRY",102,"LP1912")
  ,@("05:49:40","07:32","84_COLONIA URQUIZA-ESC 49",103,"LP1912")
  ,@("05:49:40","07:36","27_EL RETIRO",107,"LP1912")
  ,@("05:49:40","07:39","10_OLMOS",110,"LP1912")
)

$startRow = 43
$r = $startRow
foreach ($row in $newRowsLP) {
    $wsLP.Cells.Item($r,1).Value = $row[0]
    $wsLP.Cells.Item($r,2).Value = $row[1]
    $wsLP.Cells.Item($r,3).Value = $row[2]
    $wsLP.Cells.Item($r,4).Value = $row[3]
    $wsLP.Cells.Item($r,5).Value = $row[4]
    $r++
}
$lastRowLP = $r - 1

$rngLP = $wsLP.Range("A6:E$lastRowLP")
$wsLP.Sort.SortFields.Clear()
$wsLP.Sort.SortFields.Add($wsLP.Range("B6:B$lastRowLP"))
$wsLP.Sort.SetRange($rngLP)
$wsLP.Sort.Header = 0
$wsLP.Sort.Apply()

$wsLP.Range("A2").Value = "Última actualización: 05:49:40"
$wsLP.Range("A3").Value = "Total filas: 48"

# ---------------------------------------------------------------
# Hoja "LP1912-215"
# ---------------------------------------------------------------
$wsLP215 = $wb.Worksheets.Item("LP1912-215")
$wsLP215.Range("A2").Value = "Última actualización: 05:49:40"

# ---------------------------------------------------------------
# Hoja "6203-6173"
# ---------------------------------------------------------------
$wsL6 = $wb.Worksheets.Item("6203-6173")

$newRowsL6 = @(
  ,@("05:49:40","07:07","215B_LP-P MOR-1 Y 57",78,"L6173")
  ,@("05:49:40","07:35","215A_LA PLATA",106,"L6173")
)

$startRow6 = 13
$r6 = $startRow6
foreach ($row in $newRowsL6) {
    $wsL6.Cells.Item($r6,1).Value = $row[0]
    $wsL6.Cells.Item($r6,2).Value = $row[1]
    $wsL6.Cells.Item($r6,3).Value = $row[2]
    $wsL6.Cells.Item($r6,4).Value = $row[3]
    $wsL6.Cells.Item($r6,5).Value = $row[4]
    $r6++
}
$lastRowL6 = $r6 - 1

$rngL6 = $wsL6.Range("A6:E$lastRowL6")
$wsL6.Sort.SortFields.Clear()
$wsL6.Sort.SortFields.Add($wsL6.Range("B6:B$lastRowL6"))
$wsL6.Sort.SetRange($rngL6)
$wsL6.Sort.Header = 0
$wsL6.Sort.Apply()

$wsL6.Range("A2").Value = "Última actualización: 05:49:40"
$wsL6.Range("A3").Value = "Total filas: 9"
